$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.208.92'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.438.02'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +1.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9227'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '273.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3637'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3053'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.013'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06463'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9999'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.307'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.014'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001006'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('D17').Value = '1.439.41'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9425'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05647'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.342'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.248'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').Value = '20.257.73'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.019'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').Value = '1.593.40'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.77'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.052'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.762'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.07639'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7708'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.452'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05649'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.96%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.606'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.39%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.112'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9361'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01971'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.36%  '
$ws.Range('E41').Value = '  -5.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1825'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.946'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -17.88%  '
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5149'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5030'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('E49').Value = '  -3.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06351'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9875'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.97%  '
